$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.892.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.061.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.59%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -11.41%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.056.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -11.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000215"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.556.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.955.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.86%  "
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.058.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.700"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.24%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "509.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0398"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.076.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0787"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.33%  "
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.252"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.107"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +60.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₃0494"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.68%  "
